$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 10250.25
$ws.Range("I7").Value = 7333.6665
$ws.Range("K7").Value = 7333.6665
$ws.Range("M7").Value = -7221.6665
$ws.Range("H14").Value = 10250.25
$ws.Range("I14").Value = 7333.6665
$ws.Range("K14").Value = 7333.6665
$ws.Range("M14").Value = -7142.6665
$ws.Range("H138").Value = 2288.513
$ws.Range("I138").Value = 1791.625
$ws.Range("J138").Value = 2634.1738
$ws.Range("K138").Value = 5374.875
$ws.Range("L138").Value = 7902.5214
$ws.Range("M138").Value = -234.875
$ws.Range("N138").Value = -18182.5214

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30642.764
$ws.Range("I32").Value = 31279.352
$ws.Range("J32").Value = 12500
$ws.Range("K32").Value = 31279.352
$ws.Range("L32").Value = 12500
$ws.Range("M32").Value = -30992.352
$ws.Range("N32").Value = -13074
$ws.Range("H61").Value = 1923.0444
$ws.Range("I61").Value = 1777.0952
$ws.Range("K61").Value = 1777.0952
$ws.Range("M61").Value = -1565.0952
$ws.Range("H63").Value = 2924.0908
$ws.Range("I63").Value = 2401.875
$ws.Range("J63").Value = 4316.6665
$ws.Range("K63").Value = 2401.875
$ws.Range("L63").Value = 4316.6665
$ws.Range("M63").Value = -1715.875
$ws.Range("N63").Value = -5688.6665
$ws.Range("H66").Value = 2924.0908
$ws.Range("I66").Value = 2401.875
$ws.Range("J66").Value = 4316.6665
$ws.Range("K66").Value = 12009.375
$ws.Range("L66").Value = 21583.3325
$ws.Range("M66").Value = -8577.375
$ws.Range("N66").Value = -28447.3325
$ws.Range("H101").Value = 48084
$ws.Range("J101").Value = 48084
$ws.Range("L101").Value = 48084
$ws.Range("N101").Value = -54574
$ws.Range("H136").Value = 1923.0444
$ws.Range("I136").Value = 1777.0952
$ws.Range("K136").Value = 5331.2856
$ws.Range("M136").Value = -2781.2856

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 33762.5
$ws.Range("J35").Value = 33762.5
$ws.Range("L35").Value = 33762.5
$ws.Range("N35").Value = -34382.5
$ws.Range("H105").Value = 2516.3333
$ws.Range("I105").Value = 2471.25
$ws.Range("J105").Value = 2538.875
$ws.Range("K105").Value = 2471.25
$ws.Range("L105").Value = 2538.875
$ws.Range("M105").Value = -724.25
$ws.Range("N105").Value = -6032.875
$ws.Range("H107").Value = 2469.7693
$ws.Range("I107").Value = 2372.182
$ws.Range("J107").Value = 3006.5
$ws.Range("K107").Value = 2372.182
$ws.Range("L107").Value = 3006.5
$ws.Range("M107").Value = -452.1819999999998
$ws.Range("N107").Value = -6846.5
$ws.Range("H134").Value = 2591.32
$ws.Range("I134").Value = 2015.1052
$ws.Range("J134").Value = 4416
$ws.Range("K134").Value = 6045.3156
$ws.Range("L134").Value = 13248
$ws.Range("M134").Value = -3510.3156
$ws.Range("N134").Value = -18318

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4689.443
$ws.Range("I31").Value = 1846.3448
$ws.Range("J31").Value = 6700.4146
$ws.Range("K31").Value = 1846.3448
$ws.Range("L31").Value = 6700.4146
$ws.Range("M31").Value = -1551.3448
$ws.Range("N31").Value = -7290.4146
$ws.Range("H34").Value = 4689.443
$ws.Range("I34").Value = 1846.3448
$ws.Range("J34").Value = 6700.4146
$ws.Range("K34").Value = 1846.3448
$ws.Range("L34").Value = 6700.4146
$ws.Range("M34").Value = -1644.3448
$ws.Range("N34").Value = -7104.4146
$ws.Range("H58").Value = 2143.8235
$ws.Range("I58").Value = 1905.275
$ws.Range("J58").Value = 3011.2727
$ws.Range("K58").Value = 1905.275
$ws.Range("L58").Value = 3011.2727
$ws.Range("M58").Value = -1702.275
$ws.Range("N58").Value = -3417.2727
$ws.Range("H68").Value = 300000
$ws.Range("J68").Value = 300000
$ws.Range("L68").Value = 300000
$ws.Range("N68").Value = -301498
$ws.Range("H70").Value = 32036.666
$ws.Range("J70").Value = 32036.666
$ws.Range("L70").Value = 32036.666
$ws.Range("N70").Value = -32666.666
$ws.Range("H71").Value = 300000
$ws.Range("J71").Value = 300000
$ws.Range("L71").Value = 900000
$ws.Range("N71").Value = -907488
$ws.Range("H73").Value = 32036.666
$ws.Range("J73").Value = 32036.666
$ws.Range("L73").Value = 32036.666
$ws.Range("N73").Value = -34220.666
$ws.Range("H80").Value = 29536.285
$ws.Range("J80").Value = 29536.285
$ws.Range("L80").Value = 29536.285
$ws.Range("N80").Value = -31782.285
$ws.Range("H83").Value = 29536.285
$ws.Range("J83").Value = 29536.285
$ws.Range("L83").Value = 88608.855
$ws.Range("N83").Value = -99840.855
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0
$ws.Range("H132").Value = 66598.67999999999
$ws.Range("I132").Value = 2173.923
$ws.Range("J132").Value = 159656.67
$ws.Range("K132").Value = 6521.768999999999
$ws.Range("L132").Value = 478970.01
$ws.Range("M132").Value = -3991.768999999999
$ws.Range("N132").Value = -484030.01
$ws.Range("H134").Value = 3406.5925
$ws.Range("I134").Value = 1656.7894
$ws.Range("J134").Value = 7562.375
$ws.Range("K134").Value = 4970.3682
$ws.Range("L134").Value = 22687.125
$ws.Range("M134").Value = -2435.3682
$ws.Range("N134").Value = -27757.125
$ws.Range("H136").Value = 2143.8235
$ws.Range("I136").Value = 1905.275
$ws.Range("J136").Value = 3011.2727
$ws.Range("K136").Value = 5715.825000000001
$ws.Range("L136").Value = 9033.8181
$ws.Range("M136").Value = -3165.825000000001
$ws.Range("N136").Value = -14133.8181

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 250
$ws.Range("I36").Value = 250
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 750
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -581
$ws.Range("H98").Value = 2567.3333
$ws.Range("I98").Value = 2600
$ws.Range("J98").Value = 2534.6667
$ws.Range("K98").Value = 7800
$ws.Range("L98").Value = 7604.000100000001
$ws.Range("M98").Value = -6302
$ws.Range("N98").Value = -10600.0001
$ws.Range("H131").Value = 896.71
$ws.Range("I131").Value = 669.75
$ws.Range("J131").Value = 916.4457
$ws.Range("K131").Value = 2009.25
$ws.Range("L131").Value = 2749.3371
$ws.Range("M131").Value = 3030.75
$ws.Range("N131").Value = -12829.3371

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19800
$ws.Range("J15").Value = 19800
$ws.Range("L15").Value = 19800
$ws.Range("N15").Value = -20376
$ws.Range("H81").Value = 19800
$ws.Range("J81").Value = 19800
$ws.Range("L81").Value = 19800
$ws.Range("N81").Value = -21796
$ws.Range("H84").Value = 19800
$ws.Range("J84").Value = 19800
$ws.Range("L84").Value = 59400
$ws.Range("N84").Value = -69384
$ws.Range("H97").Value = 15093.667
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 15093.667
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 15093.667
$ws.Range("N97").Value = -16085.667
$ws.Range("H102").Value = 2247.375
$ws.Range("I102").Value = 2252.5
$ws.Range("K102").Value = 2252.5
$ws.Range("M102").Value = -630.5
$ws.Range("H132").Value = 2316.738
$ws.Range("I132").Value = 1871.8125
$ws.Range("J132").Value = 3740.5
$ws.Range("K132").Value = 5615.4375
$ws.Range("L132").Value = 11221.5
$ws.Range("M132").Value = -3085.4375
$ws.Range("N132").Value = -16281.5
$ws.Range("H136").Value = 66884
$ws.Range("J136").Value = 66884
$ws.Range("L136").Value = 200652
$ws.Range("N136").Value = -205752

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 6000
$ws.Range("J33").Value = 6000
$ws.Range("L33").Value = 6000
$ws.Range("N33").Value = -6580
$ws.Range("H136").Value = 2152.2
$ws.Range("I136").Value = 1873.9565
$ws.Range("J136").Value = 3066.4285
$ws.Range("K136").Value = 5621.8695
$ws.Range("L136").Value = 9199.2855
$ws.Range("M136").Value = -3071.8695
$ws.Range("N136").Value = -14299.2855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 28064.75
$ws.Range("J75").Value = 28064.75
$ws.Range("L75").Value = 28064.75
$ws.Range("N75").Value = -29936.75
$ws.Range("H78").Value = 28064.75
$ws.Range("J78").Value = 28064.75
$ws.Range("L78").Value = 84194.25
$ws.Range("N78").Value = -93554.25
$ws.Range("H104").Value = 40365.5
$ws.Range("J104").Value = 40365.5
$ws.Range("L104").Value = 40365.5
$ws.Range("N104").Value = -47353.5
$ws.Range("H123").Value = 38584
$ws.Range("J123").Value = 38584
$ws.Range("L123").Value = 38584
$ws.Range("N123").Value = -48384
$ws.Range("H132").Value = 1648
$ws.Range("I132").Value = 1123.6154
$ws.Range("J132").Value = 2216.0833
$ws.Range("K132").Value = 3370.8462
$ws.Range("L132").Value = 6648.249899999999
$ws.Range("M132").Value = -840.8462
$ws.Range("N132").Value = -11708.2499
